$d = $word.ActiveDocument

$pairs = @(
    @{old="88×62="; new="94×22="},
    @{old="81×53="; new="63×91="},
    @{old="90×91="; new="68×52="},
    @{old="49×48="; new="65×87="},
    @{old="43×62="; new="41×30="},
    @{old="15×96="; new="97×62="},
    @{old="62×60="; new="37×96="},
    @{old="39×52="; new="63×52="},
    @{old="45×69="; new="32×39="},
    @{old="44×85="; new="70×52="},
    @{old="93×88="; new="76×14="},
    @{old="20×76="; new="87×42="},
    @{old="54×72="; new="57×13="},
    @{old="74×62="; new="24×47="},
    @{old="92×88="; new="18×64="},
    @{old="53×43="; new="70×75="},
    @{old="58×82="; new="57×34="},
    @{old="19×62="; new="86×51="},
    @{old="47×47="; new="53×42="},
    @{old="79×75="; new="83×68="},
    @{old="16×31="; new="65×86="},
    @{old="32×15="; new="35×78="},
    @{old="82×98="; new="54×95="},
    @{old="85×20="; new="50×87="},
    @{old="27×37="; new="41×99="}
)

foreach ($pair in $pairs) {
    $d.Content.Find.Execute($pair.old, $true, $false, $false, $false, $false, $true, 1, $false, $pair.new, 2)
}
